# Auto-generated script applying the Ultros_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1142.5264
$ws.Range("I28").Value = 765.75
$ws.Range("K28").Value = 765.75
$ws.Range("M28").Value = -280.75
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H98").Value = 895.5278
$ws.Range("I98").Value = 892.54285
$ws.Range("K98").Value = 892.54285
$ws.Range("M98").Value = 605.45715
$ws.Range("H116").Value = 4247.1055
$ws.Range("I116").Value = 4232.9287
$ws.Range("J116").Value = 4286.8
$ws.Range("K116").Value = 4232.9287
$ws.Range("L116").Value = 4286.8
$ws.Range("M116").Value = -790.9287000000004
$ws.Range("N116").Value = -11170.8
$ws.Range("H122").Value = 895.5278
$ws.Range("I122").Value = 892.54285
$ws.Range("K122").Value = 2677.62855
$ws.Range("M122").Value = -227.6285500000004
$ws.Range("H125").Value = 22266.092
$ws.Range("I125").Value = 37899.168
$ws.Range("J125").Value = 3506.4
$ws.Range("K125").Value = 341092.512
$ws.Range("L125").Value = 31557.6
$ws.Range("M125").Value = -338632.512
$ws.Range("N125").Value = -36477.60000000001
$ws.Range("H132").Value = 12043.067
$ws.Range("I132").Value = 3027.8125
$ws.Range("J132").Value = 51382.363
$ws.Range("K132").Value = 9083.4375
$ws.Range("L132").Value = 154147.089
$ws.Range("M132").Value = -6553.4375
$ws.Range("N132").Value = -159207.089
$ws.Range("H137").Value = 3516.923
$ws.Range("I137").Value = 3643.5
$ws.Range("J137").Value = 2546.5
$ws.Range("K137").Value = 10930.5
$ws.Range("L137").Value = 7639.5
$ws.Range("M137").Value = -8380.5
$ws.Range("N137").Value = -12739.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10642424
$ws.Range("I32").Value = 11632213
$ws.Range("K32").Value = 11632213
$ws.Range("M32").Value = -11631926
$ws.Range("H45").Value = 2712.45
$ws.Range("I45").Value = 2296.6
$ws.Range("J45").Value = 3960
$ws.Range("K45").Value = 2296.6
$ws.Range("L45").Value = 3960
$ws.Range("M45").Value = -1919.6
$ws.Range("N45").Value = -4714
$ws.Range("H61").Value = 2069.5557
$ws.Range("I61").Value = 2016.8667
$ws.Range("K61").Value = 2016.8667
$ws.Range("M61").Value = -1804.8667
$ws.Range("H74").Value = 2010.6538
$ws.Range("I74").Value = 1431.5714
$ws.Range("K74").Value = 1431.5714
$ws.Range("M74").Value = -557.5714
$ws.Range("H77").Value = 2010.6538
$ws.Range("I77").Value = 1431.5714
$ws.Range("K77").Value = 7157.857
$ws.Range("M77").Value = -2789.857
$ws.Range("H102").Value = 18521130
$ws.Range("I102").Value = 2759.5334
$ws.Range("J102").Value = 111112984
$ws.Range("K102").Value = 2759.5334
$ws.Range("L102").Value = 111112984
$ws.Range("M102").Value = -1137.5334
$ws.Range("N102").Value = -111116228
$ws.Range("H132").Value = 1825.5358
$ws.Range("I132").Value = 1300.2632
$ws.Range("K132").Value = 3900.7896
$ws.Range("M132").Value = -1370.7896
$ws.Range("H136").Value = 2069.5557
$ws.Range("I136").Value = 2016.8667
$ws.Range("K136").Value = 6050.6001
$ws.Range("M136").Value = -3500.6001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 291.33334
$ws.Range("H86").Value = 3770.8386
$ws.Range("I86").Value = 2837.8462
$ws.Range("J86").Value = 4444.6665
$ws.Range("K86").Value = 2837.8462
$ws.Range("L86").Value = 4444.6665
$ws.Range("M86").Value = -1714.8462
$ws.Range("N86").Value = -6690.6665
$ws.Range("H89").Value = 3770.8386
$ws.Range("I89").Value = 2837.8462
$ws.Range("J89").Value = 4444.6665
$ws.Range("K89").Value = 14189.231
$ws.Range("L89").Value = 22223.3325
$ws.Range("M89").Value = -8573.231
$ws.Range("N89").Value = -33455.3325
$ws.Range("H134").Value = 2602.3584
$ws.Range("I134").Value = 1815.175
$ws.Range("K134").Value = 5445.525
$ws.Range("M134").Value = -2910.525

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 554.8182
$ws.Range("I15").Value = 724.75
$ws.Range("J15").Value = 101.666664
$ws.Range("K15").Value = 724.75
$ws.Range("L15").Value = 101.666664
$ws.Range("M15").Value = -554.75
$ws.Range("N15").Value = -441.666664
$ws.Range("H16").Value = 7676.25
$ws.Range("I16").Value = 7511.4
$ws.Range("J16").Value = 8500.5
$ws.Range("K16").Value = 7511.4
$ws.Range("L16").Value = 8500.5
$ws.Range("M16").Value = -7224.4
$ws.Range("N16").Value = -9074.5
$ws.Range("H31").Value = 2133.5366
$ws.Range("I31").Value = 1589.6818
$ws.Range("J31").Value = 2763.2632
$ws.Range("K31").Value = 1589.6818
$ws.Range("L31").Value = 2763.2632
$ws.Range("M31").Value = -1294.6818
$ws.Range("N31").Value = -3353.2632
$ws.Range("H34").Value = 2133.5366
$ws.Range("I34").Value = 1589.6818
$ws.Range("J34").Value = 2763.2632
$ws.Range("K34").Value = 1589.6818
$ws.Range("L34").Value = 2763.2632
$ws.Range("M34").Value = -1387.6818
$ws.Range("N34").Value = -3167.2632
$ws.Range("H58").Value = 1365.4849
$ws.Range("I58").Value = 1115.3667
$ws.Range("K58").Value = 1115.3667
$ws.Range("M58").Value = -912.3667
$ws.Range("H86").Value = 75828.836
$ws.Range("I86").Value = 71600.8
$ws.Range("K86").Value = 71600.8
$ws.Range("M86").Value = -70477.8
$ws.Range("H89").Value = 75828.836
$ws.Range("I89").Value = 71600.8
$ws.Range("K89").Value = 358004
$ws.Range("M89").Value = -352388
$ws.Range("H113").Value = 7676.25
$ws.Range("I113").Value = 7511.4
$ws.Range("J113").Value = 8500.5
$ws.Range("K113").Value = 7511.4
$ws.Range("L113").Value = 8500.5
$ws.Range("M113").Value = -5341.4
$ws.Range("N113").Value = -12840.5
$ws.Range("H134").Value = 2630.617
$ws.Range("I134").Value = 2442.805
$ws.Range("J134").Value = 3914
$ws.Range("K134").Value = 7328.414999999999
$ws.Range("L134").Value = 11742
$ws.Range("M134").Value = -4793.414999999999
$ws.Range("N134").Value = -16812
$ws.Range("H136").Value = 1365.4849
$ws.Range("I136").Value = 1115.3667
$ws.Range("K136").Value = 3346.1001
$ws.Range("M136").Value = -796.1001000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4444.3335
$ws.Range("I80").Value = 3999.6667
$ws.Range("J80").Value = 4666.6665
$ws.Range("K80").Value = 11999.0001
$ws.Range("L80").Value = 13999.9995
$ws.Range("M80").Value = -11063.0001
$ws.Range("N80").Value = -15871.9995
$ws.Range("H83").Value = 4444.3335
$ws.Range("I83").Value = 3999.6667
$ws.Range("J83").Value = 4666.6665
$ws.Range("K83").Value = 35997.0003
$ws.Range("L83").Value = 41999.9985
$ws.Range("M83").Value = -31317.0003
$ws.Range("N83").Value = -51359.9985
$ws.Range("H122").Value = 725.9231
$ws.Range("I122").Value = 769.5714
$ws.Range("J122").Value = 675
$ws.Range("K122").Value = 6926.1426
$ws.Range("L122").Value = 6075
$ws.Range("M122").Value = -4476.1426
$ws.Range("N122").Value = -10975
$ws.Range("H132").Value = 2333
$ws.Range("I132").Value = 2333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 20997
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -18467
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value = 7274.5
$ws.Range("I138").Value = 6026
$ws.Range("J138").Value = 7898.75
$ws.Range("K138").Value = 18078
$ws.Range("L138").Value = 23696.25
$ws.Range("M138").Value = -12938
$ws.Range("N138").Value = -33976.25
$ws.Range("H140").Value = 16667789
$ws.Range("I140").Value = 16667789
$ws.Range("K140").Value = 50003367
$ws.Range("M140").Value = -49998187

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3189.0571
$ws.Range("I132").Value = 2062.4717
$ws.Range("K132").Value = 6187.4151
$ws.Range("M132").Value = -3657.4151
$ws.Range("H136").Value = 4937.857
$ws.Range("I136").Value = 5128.5884
$ws.Range("J136").Value = 4127.25
$ws.Range("K136").Value = 15385.7652
$ws.Range("L136").Value = 12381.75
$ws.Range("M136").Value = -12835.7652
$ws.Range("N136").Value = -17481.75
$ws.Range("H141").Value = 134166.75
$ws.Range("J141").Value = 138181.9
$ws.Range("L141").Value = 138181.9
$ws.Range("N141").Value = -148541.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H136").Value = 2674.4285
$ws.Range("I136").Value = 1911.8334
$ws.Range("K136").Value = 5735.5002
$ws.Range("M136").Value = -3185.5002
